$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the "teclas acceso rapido" task (row 2) as "en proceso" in column C,
# matching the existing pattern used for other in-progress tasks (e.g. row 5)
$ws.Range("C2").Value = "en proceso"

# Update the selection/cursor position left by the user after the edit
$ws.Range("A3").Select()
